$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 303; $r -le 387; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 * 0.55
}
